# Updated symbol list on Thu Dec 22 02:11:21 UTC 2022 with GitHub Actions
#
# This script updates the "Price" (column D), "Volume(1h)" (column E, only
# two special rows) and "Hora" (column G) columns of the crypto price
# table on the active sheet, matching a refreshed data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values, keyed by row number. All of these cells are
# stored as text in the workbook, so we force a text number format before
# assigning the value (and clear the format afterwards) to avoid Excel
# auto-converting the numeric-looking strings into real numbers.
$priceUpdates = @(
    @{ Row = 2;  Value = "246.79" },
    @{ Row = 3;  Value = "22.58" },
    @{ Row = 4;  Value = "5.335" },
    @{ Row = 5;  Value = "0.05722" },
    @{ Row = 6;  Value = "3.428" },
    @{ Row = 7;  Value = "0.8080" },
    @{ Row = 8;  Value = "0.8646" },
    @{ Row = 9;  Value = "0.1424" },
    @{ Row = 10; Value = "0.07366" },
    @{ Row = 11; Value = "0.03042" },
    @{ Row = 12; Value = "0.03116" },
    @{ Row = 13; Value = "0.09386" },
    @{ Row = 14; Value = "3.899" },
    @{ Row = 15; Value = "0.001582" },
    @{ Row = 16; Value = "0.04795" },
    @{ Row = 17; Value = "0.0005854" },
    @{ Row = 18; Value = "0.006154" },
    @{ Row = 19; Value = "0.005161" },
    @{ Row = 20; Value = "0.0009964" },
    @{ Row = 21; Value = "0.0001500" },
    @{ Row = 22; Value = "3.701" },
    @{ Row = 23; Value = "6.313" },
    @{ Row = 24; Value = "2.201" },
    @{ Row = 25; Value = "0.3264" },
    @{ Row = 26; Value = "0.1295" },
    @{ Row = 40; Value = "0.03944" },
    @{ Row = 41; Value = "0.006776" },
    @{ Row = 42; Value = "0.1068" },
    @{ Row = 43; Value = "0.002391" },
    @{ Row = 44; Value = "0.007991" },
    @{ Row = 45; Value = "0.00005614" },
    @{ Row = 47; Value = "0.3602" },
    @{ Row = 48; Value = "0.1808" },
    @{ Row = 49; Value = "0.00002101" },
    @{ Row = 50; Value = "0.01011" }
)

foreach ($update in $priceUpdates) {
    $cell = $ws.Cells.Item($update.Row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.ClearFormats()
}

# Two "Volume(1h)" (column E) labels had their "Worstin24h" / "Bestin24h"
# suffix move from one coin to another.
$e17 = $ws.Cells.Item(17, 5)
$e17.NumberFormat = "@"
$e17.Value = "16OneONE"
$e17.ClearFormats()

$e47 = $ws.Cells.Item(47, 5)
$e47.NumberFormat = "@"
$e47.Value = "46CoinbaseStockTokenCOINWorstin24h"
$e47.ClearFormats()

# Every data row's "Hora" (column G) value moves from "0" to "2".
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$gRange.Value = "2"
$gRange.ClearFormats()
